$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("D2", "30.402.22"),
    @("E2", "  -1.92%  "),
    @("D3", "1.905.71"),
    @("E3", "  -2.66%  "),
    @("D4", "'0.9995"),
    @("E4", "  -0.15%  "),
    @("D5", "'238.28"),
    @("E5", "  -2.43%  "),
    @("D6", "'1.000"),
    @("E6", "  -0.04%  "),
    @("D7", "'0.4729"),
    @("E7", "  -2.58%  "),
    @("D8", "'0.2826"),
    @("E8", "  -3.78%  "),
    @("D9", "'0.06644"),
    @("E9", "  -5.07%  "),
    @("D10", "'18.59"),
    @("E10", "  -5.70%  "),
    @("D11", "'100.42"),
    @("E11", "  -6.21%  "),
    @("D12", "'0.07717"),
    @("E12", "  -0.84%  "),
    @("D13", "1.902.38"),
    @("E13", "  -2.73%  "),
    @("D14", "'5.181"),
    @("E14", "  -4.31%  "),
    @("D15", "'0.6660"),
    @("E15", "  -5.49%  "),
    @("D16", "30.402.62"),
    @("E16", "  -1.94%  "),
    @("D17", "'254.56"),
    @("E17", "  -8.70%  "),
    @("E18", "  +0.00%  "),
    @("D19", "'0.000007423"),
    @("E19", "  -4.73%  "),
    @("D20", "'12.58"),
    @("E20", "  -5.36%  "),
    @("D21", "'5.357"),
    @("E21", "  -3.83%  "),
    @("D22", "'0.9991"),
    @("E22", "  -0.22%  "),
    @("B23", "BitDAO"),
    @("C23", "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"),
    @("D23", "'0.4560"),
    @("E23", "  -8.75%  "),
    @("B24", "Chainlink"),
    @("C24", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"),
    @("D24", "'6.274"),
    @("E24", "  -3.87%  "),
    @("B25", "Cosmos"),
    @("C25", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"),
    @("D25", "'9.321"),
    @("E25", "  -4.85%  "),
    @("B26", "Monero"),
    @("C26", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"),
    @("D26", "'166.77"),
    @("E26", "  -1.53%  "),
    @("B27", "EthereumClassic"),
    @("C27", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"),
    @("D27", "'18.86"),
    @("E27", "  -4.59%  "),
    @("B28", "LidoDAOToken"),
    @("C28", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"),
    @("D28", "'2.042"),
    @("E28", "  -6.66%  "),
    @("D29", "'0.1010"),
    @("E29", "  -3.89%  "),
    @("D30", "'4.705"),
    @("E30", "  +1.64%  "),
    @("B31", "Toncoin"),
    @("C31", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"),
    @("D31", "'1.376"),
    @("E31", "  -0.74%  "),
    @("B32", "PancakeSwap"),
    @("C32", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"),
    @("D32", "'1.506"),
    @("E32", "  -4.06%  "),
    @("B33", "InternetComputer(DFINITY)"),
    @("C33", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"),
    @("D33", "'4.245"),
    @("E33", "  -4.46%  "),
    @("B34", "Hedera"),
    @("C34", "https://coinranking.com/coin/jad286TjB+hedera-hbar"),
    @("D34", "'0.04706"),
    @("E34", "  -4.20%  "),
    @("B35", "ImmutableX"),
    @("C35", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"),
    @("D35", "'0.7232"),
    @("E35", "  -3.52%  "),
    @("B36", "ARBITRUM"),
    @("C36", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"),
    @("D36", "'1.102"),
    @("E36", "  -5.82%  "),
    @("B37", "Frax"),
    @("C37", "https://coinranking.com/coin/KfWtaeV1W+frax-frax"),
    @("D37", "'0.9996"),
    @("E37", "  -0.09%  "),
    @("B38", "HuobiToken"),
    @("C38", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"),
    @("D38", "'2.707"),
    @("E38", "  -1.07%  "),
    @("B39", "VeChain"),
    @("C39", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"),
    @("D39", "'0.01906"),
    @("E39", "  -4.99%  "),
    @("B40", "MXToken"),
    @("C40", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"),
    @("D40", "'2.594"),
    @("E40", "  -3.58%  "),
    @("B41", "FraxShare"),
    @("C41", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"),
    @("D41", "'6.178"),
    @("E41", "  -5.82%  "),
    @("B42", "Aave"),
    @("C42", "https://coinranking.com/coin/ixgUfzmLR+aave-aave"),
    @("D42", "'72.35"),
    @("E42", "  -7.32%  "),
    @("B43", "RenderToken"),
    @("C43", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"),
    @("D43", "'1.959"),
    @("E43", "  -8.28%  "),
    @("B44", "Quant"),
    @("C44", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"),
    @("D44", "'105.83"),
    @("E44", "  -3.34%  "),
    @("B45", "TrustWalletToken"),
    @("C45", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"),
    @("D45", "'0.8544"),
    @("E45", "  -4.73%  "),
    @("B46", "PaxDollar"),
    @("C46", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"),
    @("D46", "'1.001"),
    @("E46", "  +0.00%  "),
    @("B47", "TheSandbox"),
    @("C47", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"),
    @("D47", "'0.4211"),
    @("E47", "  -5.44%  "),
    @("B48", "Maker"),
    @("C48", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"),
    @("D48", "'988.15"),
    @("E48", "  -0.30%  "),
    @("B49", "Aptos"),
    @("C49", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"),
    @("D49", "'7.356"),
    @("E49", "  -7.68%  "),
    @("B50", "Algorand"),
    @("C50", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"),
    @("D50", "'0.1190"),
    @("E50", "  -4.65%  "),
    @("B51", "Elrond"),
    @("C51", "https://coinranking.com/coin/omwkOTglq+elrond-egld"),
    @("D51", "'34.35"),
    @("E51", "  -4.42%  "),
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}
